$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in M1: "SAVE_INTV" -> "INTV_PRINT"
$ws.Range("M1").Value = "INTV_PRINT"

# Update L2 value: 22.6 -> 30
$ws.Range("L2").Value = 30

# Update N2 formula: multiplier 2 -> 11
$ws.Range("N2").Formula = "=D2*L2*11/M2/1000"

# Update the active cell selection to M2
$ws.Range("M2").Select()
